$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set all new / changed text values first, in the same order the author
# --- typed them, so new shared-string entries land in the matching order. ---

# Row 18 "Save button": trim trailing space off "Save ".
$ws.Range("B18").Value = "Save"

# New row 19: Select Role / Admin / ESS
$ws.Range("A19").Value = "Select Role"
$ws.Range("B19").Value = "Admin"
$ws.Range("C19").Value = "ESS"

# New row 20: Seletc status / Enabled / Disabled (author filled C before B)
$ws.Range("A20").Value = "Seletc status"
$ws.Range("C20").Value = "Disabled"
$ws.Range("B20").Value = "Enabled"

# Row 10 "Employee Name": swap the sample employee name (edited last).
$ws.Range("B10").Value = "Lisa Michelle Walker"

# --- Now bring formatting in line: the label/value cells in this sheet use
# --- a Consolas-font style (same as most other "B" value cells). ---

$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("C20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Scroll the view down a bit (sheet view's topLeftCell moves from A4 to A7)
# while keeping the original active-cell selection at B13.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B13").Select() | Out-Null
